$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'275.95"
$ws.Range("E2").Formula = "'-0.63%"
$ws.Range("D3").Formula = "'27.35"
$ws.Range("E3").Formula = "'2.02%"
$ws.Range("D4").Formula = "'4.795"
$ws.Range("E4").Formula = "'-2.54%"
$ws.Range("D5").Formula = "'0.06357"
$ws.Range("E5").Formula = "'-0.50%"
$ws.Range("D6").Formula = "'6.952"
$ws.Range("E6").Formula = "'-0.25%"
$ws.Range("D7").Formula = "'1.310"
$ws.Range("E7").Formula = "'2.45%"
$ws.Range("D8").Formula = "'0.8772"
$ws.Range("E8").Formula = "'-0.95%"
$ws.Range("D9").Formula = "'0.1549"
$ws.Range("E9").Formula = "'4.77%"
$ws.Range("D10").Formula = "'0.05060"
$ws.Range("E10").Formula = "'-2.49%"
$ws.Range("D11").Formula = "'0.07595"
$ws.Range("E11").Formula = "'2.30%"
$ws.Range("D12").Formula = "'0.03024"
$ws.Range("E12").Formula = "'-4.00%"
$ws.Range("D13").Formula = "'0.09041"
$ws.Range("E13").Formula = "'-0.26%"
$ws.Range("D14").Formula = "'0.001573"
$ws.Range("E14").Formula = "'0.74%"
$ws.Range("D15").Formula = "'0.0006408"
$ws.Range("E15").Formula = "'1.10%"
$ws.Range("D16").Formula = "'0.005801"
$ws.Range("E16").Formula = "'-3.52%"
$ws.Range("D17").Formula = "'3.454"
$ws.Range("E17").Formula = "'-1.02%"
$ws.Range("D18").Formula = "'3.302"
$ws.Range("E18").Formula = "'-1.53%"
$ws.Range("E19").Formula = "'-0.39%"
$ws.Range("E20").Formula = "'-1.20%"
$ws.Range("D21").Formula = "'0.1336"
$ws.Range("E21").Formula = "'0.38%"
$ws.Range("D22").Formula = "'3.954"
$ws.Range("E22").Formula = "'1.22%"
$ws.Range("E23").Formula = "'1.33%"
$ws.Range("D24").Formula = "'0.001172"
$ws.Range("E24").Formula = "'-0.84%"
$ws.Range("D25").Formula = "'0.003865"
$ws.Range("E25").Formula = "'5.11%"
$ws.Range("D26").Formula = "'0.0001200"
$ws.Range("E26").Formula = "'-0.18%"
$ws.Range("E27").Formula = "'19.67%"
$ws.Range("D40").Formula = "'0.04160"
$ws.Range("E40").Formula = "'2.29%"
$ws.Range("D41").Formula = "'0.006852"
$ws.Range("E41").Formula = "'2.79%"
$ws.Range("E42").Formula = "'0.33%"
$ws.Range("D43").Formula = "'0.002020"
$ws.Range("E43").Formula = "'-14.57%"
$ws.Range("D44").Formula = "'0.01118"
$ws.Range("E44").Formula = "'-12.87%"
$ws.Range("D45").Formula = "'0.00005157"
$ws.Range("E45").Formula = "'-2.06%"
$ws.Range("D47").Formula = "'0.02300"
$ws.Range("E47").Formula = "'8.29%"
